$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 74; $r += 6) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2025/12/13") {
        # Force the new value to be stored as text (matching the original
        # inline-string cell) instead of letting Excel auto-parse the
        # date-shaped text into a real date serial number.
        $cell.NumberFormat = "@"
        $cell.Value = "2025/12/14"
        $cell.NumberFormat = "General"
        $cell.Style = "Normal"
    }
}
